$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Highlight the existing "Descricao" / "Descricao detalhada" rows
#    (B6, B7) in yellow - these already have the plain border style,
#    we just add the fill on top. This mints fill index 2 (yellow)
#    and cellXf index 3 (fillId=2, borderId=1).
# ------------------------------------------------------------------
$ws.Range("B6").Interior.Color = 65535
$ws.Range("B7").Interior.Color = 65535

# ------------------------------------------------------------------
# 2) New "Tabela descricao" mini-table in column D (rows 14-17)
#    D14/D15 = yellow header block with border on left/top/bottom
#    only (no right edge) -> mints border index 2 and cellXf index 4.
#    D16/D17 = yellow with the regular full border (reuses border 1
#    and fill 2 -> cellXf index 3, same shape as B6/B7 above).
# ------------------------------------------------------------------

# D16 / D17 first so fill/border caches are populated in the same
# order the source workbook used (keeps shared-string order aligned).
$ws.Range("B4").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = "Tabela descricao"
$ws.Range("D14").Interior.Color = 65535
$ws.Range("D14").Borders.Item(10).LineStyle = -4142

$ws.Range("B4").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = "id"
$ws.Range("D15").Interior.Color = 65535
$ws.Range("D15").Borders.Item(10).LineStyle = -4142

$ws.Range("B4").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "descicao-resumida"
$ws.Range("D16").Interior.Color = 65535

$ws.Range("B4").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "descicao-detalhada"
$ws.Range("D17").Interior.Color = 65535

# ------------------------------------------------------------------
# 3) New "perfil" mini-table in column E (rows 13-15)
# ------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("E13").Value = "perfil"

$ws.Range("B4").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").Value = "idperfil"

$ws.Range("B4").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = "nomeperfil"

# ------------------------------------------------------------------
# 4) New "categoria" mini-table in column F (rows 13-15)
# ------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("F13").Value = "categoria"

$ws.Range("B4").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Value = "idcategoria"

$ws.Range("B4").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = "nomecategoria"

# ------------------------------------------------------------------
# 5) Scroll / selection matches the saved view in the edited workbook
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("E20").Select()
